# Natmi following Dr Hou advice
#
# The LR-pairs table for Mmp9 -> Cd44 gains "ECs" as an additional sending
# cluster (it was already a target cluster), so the data block grows from
# 3 senders x 4 targets (12 rows) to 4 senders x 4 targets (16 rows), and
# all of the per-row statistics are recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Mmp9", "Cd44", "ECs", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 239.0839323333333, 717.251797, 0.4086975387666237, 0.4086975387666237, 42719.41310750408, 384474.7179675367, 0.3979062503196, 0.3979062503196),
    @("ECs", "Mmp9", "Cd44", "FAPs", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 117.0512696666667, 351.153809, 0.2000910950200451, 0.2000910950200451, 20914.66998575479, 188232.0298717931, 0.1948078708328911, 0.1948078708328911),
    @("ECs", "Mmp9", "Cd44", "M2", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 171.15883, 513.47649, 0.2925842480357353, 0.2925842480357353, 30582.58534736189, 275243.268126257, 0.2848588258931467, 0.2848588258931467),
    @("ECs", "Mmp9", "Cd44", "sCs", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 57.695868, 173.087604, 0.09862711817759588, 0.09862711817759588, 10309.07261576938, 92781.65354192437, 0.09602295842619382, 0.09602295842619382),
    @("FAPs", "Mmp9", "Cd44", "ECs", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 239.0839323333333, 717.251797, 0.4086975387666237, 0.4086975387666237, 389.031636678424, 3501.284730105816, 0.003623601275065665, 0.003623601275065665),
    @("FAPs", "Mmp9", "Cd44", "FAPs", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 117.0512696666667, 351.153809, 0.2000910950200451, 0.2000910950200451, 190.463016771128, 1714.167150940152, 0.001774051170535534, 0.001774051170535534),
    @("FAPs", "Mmp9", "Cd44", "M2", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 171.15883, 513.47649, 0.2925842480357353, 0.2925842480357353, 278.50554036408, 2506.54986327672, 0.002594115583484892, 0.002594115583484892),
    @("FAPs", "Mmp9", "Cd44", "sCs", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 57.695868, 173.087604, 0.09862711817759588, 0.09862711817759588, 93.881331708768, 844.931985378912, 0.0008744494822819677, 0.0008744494822819677),
    @("M2", "Mmp9", "Cd44", "ECs", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 239.0839323333333, 717.251797, 0.4086975387666237, 0.4086975387666237, 618.5380260717723, 5566.84223464595, 0.005761318537193866, 0.005761318537193867),
    @("M2", "Mmp9", "Cd44", "FAPs", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 117.0512696666667, 351.153809, 0.2000910950200451, 0.2000910950200451, 302.8252906091278, 2725.42761548215, 0.002820639777634373, 0.002820639777634374),
    @("M2", "Mmp9", "Cd44", "M2", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 171.15883, 513.47649, 0.2925842480357353, 0.2925842480357353, 442.8078617401667, 3985.2707556615, 0.004124495236712863, 0.004124495236712863),
    @("M2", "Mmp9", "Cd44", "sCs", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 57.695868, 173.087604, 0.09862711817759588, 0.09862711817759588, 149.2659417006, 1343.3934753054, 0.001390324605186972, 0.001390324605186972),
    @("sCs", "Mmp9", "Cd44", "ECs", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 239.0839323333333, 717.251797, 0.4086975387666237, 0.4086975387666237, 150.9884367025377, 1358.895930322839, 0.001406368634764122, 0.001406368634764122),
    @("sCs", "Mmp9", "Cd44", "FAPs", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 117.0512696666667, 351.153809, 0.2000910950200451, 0.2000910950200451, 73.92127128132033, 665.291441531883, 0.0006885332389840651, 0.0006885332389840652),
    @("sCs", "Mmp9", "Cd44", "M2", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 171.15883, 513.47649, 0.2925842480357353, 0.2925842480357353, 108.09176475107, 972.8258827596301, 0.001006811322390836, 0.001006811322390836),
    @("sCs", "Mmp9", "Cd44", "sCs", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 57.695868, 173.087604, 0.09862711817759588, 0.09862711817759588, 36.436613822172, 327.929524399548, 0.0003393856639331263, 0.0003393856639331263)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
